$d = $word.ActiveDocument

# wdReplace constants: 0 = wdReplaceNone, 1 = wdReplaceOne, 2 = wdReplaceAll

# 1) "2. Количество преподавателей-руководителей практики" : 434 -> 456
$d.Content.Find.Execute("434", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "456", 1)

# 2) "6. Количество студентов, сдавших зачет на «хорошо» и «отлично»" : 34 -> 456
$d.Content.Find.Execute("34", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "456", 1)

# 3) & 4) "12. Количество лекций, прочитанных для студентов" and
#         "- в том числе преподавателями института" : 323 -> 232 (both occurrences)
$d.Content.Find.Execute("323", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "232", 2)
